# Update LR-pair data rows (2-7 changed, 8-10 newly added) per Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ntn1"
$ws.Cells.Item(2,3).Value = "Adora2b"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 2.082653666666667
$ws.Cells.Item(2,8).Value = 6.247961
$ws.Cells.Item(2,9).Value = 0.0472190032704503
$ws.Cells.Item(2,10).Value = 0.0472190032704503
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.3622106666666667
$ws.Cells.Item(2,14).Value = 1.086632
$ws.Cells.Item(2,15).Value = 0.06114405336448104
$ws.Cells.Item(2,16).Value = 0.06114405336448104
$ws.Cells.Item(2,17).Value = 0.7543593730391112
$ws.Cells.Item(2,18).Value = 6.789234357352
$ws.Cells.Item(2,19).Value = 0.002887161255786018
$ws.Cells.Item(2,20).Value = 0.002887161255786018

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ntn1"
$ws.Cells.Item(3,3).Value = "Adora2b"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.082653666666667
$ws.Cells.Item(3,8).Value = 6.247961
$ws.Cells.Item(3,9).Value = 0.0472190032704503
$ws.Cells.Item(3,10).Value = 0.0472190032704503
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.875117333333333
$ws.Cells.Item(3,14).Value = 5.625351999999999
$ws.Cells.Item(3,15).Value = 0.316534781675848
$ws.Cells.Item(3,16).Value = 0.316534781675848
$ws.Cells.Item(3,17).Value = 3.905219989696889
$ws.Cells.Item(3,18).Value = 35.146979907272
$ws.Cells.Item(3,19).Value = 0.01494645689116314
$ws.Cells.Item(3,20).Value = 0.01494645689116314

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ntn1"
$ws.Cells.Item(4,3).Value = "Adora2b"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.082653666666667
$ws.Cells.Item(4,8).Value = 6.247961
$ws.Cells.Item(4,9).Value = 0.0472190032704503
$ws.Cells.Item(4,10).Value = 0.0472190032704503
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 3.686562333333333
$ws.Cells.Item(4,14).Value = 11.059687
$ws.Cells.Item(4,15).Value = 0.622321164959671
$ws.Cells.Item(4,16).Value = 0.622321164959671
$ws.Cells.Item(4,17).Value = 7.67783256091189
$ws.Cells.Item(4,18).Value = 69.100493048207
$ws.Cells.Item(4,19).Value = 0.02938538512350115
$ws.Cells.Item(4,20).Value = 0.02938538512350115

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ntn1"
$ws.Cells.Item(5,3).Value = "Adora2b"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 37.28222
$ws.Cells.Item(5,8).Value = 111.84666
$ws.Cells.Item(5,9).Value = 0.8452818134314446
$ws.Cells.Item(5,10).Value = 0.8452818134314446
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.3622106666666667
$ws.Cells.Item(5,14).Value = 1.086632
$ws.Cells.Item(5,15).Value = 0.06114405336448104
$ws.Cells.Item(5,16).Value = 0.06114405336448104
$ws.Cells.Item(5,17).Value = 13.50401776101333
$ws.Cells.Item(5,18).Value = 121.53615984912
$ws.Cells.Item(5,19).Value = 0.05168395630847755
$ws.Cells.Item(5,20).Value = 0.05168395630847755

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Ntn1"
$ws.Cells.Item(6,3).Value = "Adora2b"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 37.28222
$ws.Cells.Item(6,8).Value = 111.84666
$ws.Cells.Item(6,9).Value = 0.8452818134314446
$ws.Cells.Item(6,10).Value = 0.8452818134314446
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.875117333333333
$ws.Cells.Item(6,14).Value = 5.625351999999999
$ws.Cells.Item(6,15).Value = 0.316534781675848
$ws.Cells.Item(6,16).Value = 0.316534781675848
$ws.Cells.Item(6,17).Value = 69.90853694714664
$ws.Cells.Item(6,18).Value = 629.1768325243198
$ws.Cells.Item(6,19).Value = 0.2675610942690872
$ws.Cells.Item(6,20).Value = 0.2675610942690872

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ntn1"
$ws.Cells.Item(7,3).Value = "Adora2b"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 37.28222
$ws.Cells.Item(7,8).Value = 111.84666
$ws.Cells.Item(7,9).Value = 0.8452818134314446
$ws.Cells.Item(7,10).Value = 0.8452818134314446
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.686562333333333
$ws.Cells.Item(7,14).Value = 11.059687
$ws.Cells.Item(7,15).Value = 0.622321164959671
$ws.Cells.Item(7,16).Value = 0.622321164959671
$ws.Cells.Item(7,17).Value = 137.4432279550466
$ws.Cells.Item(7,18).Value = 1236.98905159542
$ws.Cells.Item(7,19).Value = 0.5260367628538799
$ws.Cells.Item(7,20).Value = 0.5260367628538799

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Ntn1"
$ws.Cells.Item(8,3).Value = "Adora2b"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 4.741387
$ws.Cells.Item(8,8).Value = 14.224161
$ws.Cells.Item(8,9).Value = 0.107499183298105
$ws.Cells.Item(8,10).Value = 0.107499183298105
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.3622106666666667
$ws.Cells.Item(8,14).Value = 1.086632
$ws.Cells.Item(8,15).Value = 0.06114405336448104
$ws.Cells.Item(8,16).Value = 0.06114405336448104
$ws.Cells.Item(8,17).Value = 1.717380946194667
$ws.Cells.Item(8,18).Value = 15.456428515752
$ws.Cells.Item(8,19).Value = 0.006572935800217464
$ws.Cells.Item(8,20).Value = 0.006572935800217463

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Ntn1"
$ws.Cells.Item(9,3).Value = "Adora2b"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 4.741387
$ws.Cells.Item(9,8).Value = 14.224161
$ws.Cells.Item(9,9).Value = 0.107499183298105
$ws.Cells.Item(9,10).Value = 0.107499183298105
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.875117333333333
$ws.Cells.Item(9,14).Value = 5.625351999999999
$ws.Cells.Item(9,15).Value = 0.316534781675848
$ws.Cells.Item(9,16).Value = 0.316534781675848
$ws.Cells.Item(9,17).Value = 8.890656947741332
$ws.Cells.Item(9,18).Value = 80.01591252967199
$ws.Cells.Item(9,19).Value = 0.03402723051559765
$ws.Cells.Item(9,20).Value = 0.03402723051559765

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Ntn1"
$ws.Cells.Item(10,3).Value = "Adora2b"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 4.741387
$ws.Cells.Item(10,8).Value = 14.224161
$ws.Cells.Item(10,9).Value = 0.107499183298105
$ws.Cells.Item(10,10).Value = 0.107499183298105
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.686562333333333
$ws.Cells.Item(10,14).Value = 11.059687
$ws.Cells.Item(10,15).Value = 0.622321164959671
$ws.Cells.Item(10,16).Value = 0.622321164959671
$ws.Cells.Item(10,17).Value = 17.47941872195634
$ws.Cells.Item(10,18).Value = 157.314768497607
$ws.Cells.Item(10,19).Value = 0.06689901698228994
$ws.Cells.Item(10,20).Value = 0.06689901698228994

